$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = "Wed Nov 01 15:54:40 EDT 2023"
$ws.Range("B3").Value = "Wed Nov 01 15:54:52 EDT 2023"
$ws.Range("B4").Value = "Wed Nov 01 15:55:05 EDT 2023"
